$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 662.1
$ws.Range("C2").Value = 673.25
$ws.Range("B3").Value = 3317.15
$ws.Range("C3").Value = 3318.75
$ws.Range("B4").Value = 535.1
$ws.Range("C4").Value = 529.8
$ws.Range("B5").Value = 1619.3
$ws.Range("C5").Value = 1613.4
$ws.Range("B6").Value = 6602.45
$ws.Range("C6").Value = 6571.45
$ws.Range("B7").Value = 199
$ws.Range("C7").Value = 199.85
$ws.Range("B8").Value = 273.2
$ws.Range("C8").Value = 271.15
$ws.Range("B9").Value = 47691.25
$ws.Range("C9").Value = 47578
$ws.Range("B10").Value = 931.8
$ws.Range("C10").Value = 909.9
$ws.Range("B11").Value = 3751.45
$ws.Range("C11").Value = 3843.55
$ws.Range("B12").Value = 153.6
$ws.Range("C12").Value = 152.6
$ws.Range("B13").Value = 1637.95
$ws.Range("C13").Value = 1644.1
$ws.Range("B14").Value = 524.95
$ws.Range("C14").Value = 519.1
$ws.Range("B15").Value = 1531.3
$ws.Range("C15").Value = 1526.65
$ws.Range("B16").Value = 839.45
$ws.Range("C16").Value = 822.85
$ws.Range("B17").Value = 657.75
$ws.Range("C17").Value = 648.75
$ws.Range("B18").Value = 1921.05
$ws.Range("C18").Value = 1972.95
$ws.Range("B19").Value = 287.5
$ws.Range("C19").Value = 287.75
$ws.Range("B20").Value = 22497.95
$ws.Range("C20").Value = 22443.2
$ws.Range("B21").Value = 353.85
$ws.Range("C21").Value = 343.5
$ws.Range("B22").Value = 772.05
$ws.Range("C22").Value = 769.15
$ws.Range("B23").Value = 639.05
$ws.Range("C23").Value = 613.9
$ws.Range("B24").Value = 987.2
$ws.Range("C24").Value = 977.4
$ws.Range("B25").Value = 391.75
$ws.Range("C25").Value = 377.4
$ws.Range("B26").Value = 153.1
$ws.Range("C26").Value = 149.95

# Apply number format style (s="1") to newly formatted B cells
$ws.Range("B3").NumberFormat = "#,##0.00"
$ws.Range("B6").NumberFormat = "#,##0.00"
$ws.Range("B9").NumberFormat = "#,##0.00"
$ws.Range("B11").NumberFormat = "#,##0.00"
$ws.Range("B13").NumberFormat = "#,##0.00"
$ws.Range("B15").NumberFormat = "#,##0.00"
$ws.Range("B18").NumberFormat = "#,##0.00"
$ws.Range("B20").NumberFormat = "#,##0.00"

# Update selection
$ws.Range("K11:L12").Select()
